$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert 3 new rows above the old row 65 ("considering the bot is
# open-source ..." block), shifting the rest of the checklist down by 3.
# ---------------------------------------------------------------------------
$ws.Range("A65:A67").EntireRow.Insert() | Out-Null

# New row 65: "Shift google dot by one every day at…. Urrrr 4:30am?"
# Cell gets the same grey "section header" look already used on A64
# ("Queries"), so copy that formatting across before writing the text.
$ws.Range("A64").Copy() | Out-Null
$ws.Range("A65").PasteSpecial(-4122) | Out-Null
$ws.Range("A65").Value = "Shift google dot by one every day at…. Urrrr 4:30am?"

# New row 66: "Ranks are held in Disc already" (bold)
$ws.Range("A66").Value = "Ranks are held in Disc already"
$ws.Range("A66").Font.Bold = $true

# New row 67: "ask offs about visibility"
$ws.Range("A67").Value = "ask offs about visibility"

# Row 68 keeps its original text (shifted down from old row 65) and gains
# two new cells: a bold "TODO:" label in column C and a note in column D.
$ws.Range("C68").Value = "TODO:"
$ws.Range("C68").Font.Bold = $true
$ws.Range("D68").Value = "get spreadsheet rolling by date"

# New final bullet appended after the old last row (now row 94).
$ws.Range("A95").Value = "SQL injection protection"

# ---------------------------------------------------------------------------
# Column widths: columns D & E widen to a uniform 30 characters.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 29.166666666666668
$ws.Columns.Item(5).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# Selection / view state
# ---------------------------------------------------------------------------
$ws.Range("D91").Select() | Out-Null
